# The commit swaps the presentation's theme: the slide-master theme
# (ppt/theme/theme1.xml, currently the custom "Integral" / "Red Violet"
# colour scheme) is replaced with the stock Office "Office Theme" colour
# scheme (the colours that used to live in ppt/theme/theme2.xml, the
# notes-master theme). Font scheme / format scheme are identical between
# the two themes already, so only the 12 theme colours need to change.

function HexToRgbVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" colour scheme, in PowerPoint ThemeColorScheme order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRgbVal $officeColors[$i - 1]
}
